$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-09-25"

# Update the header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 09-25)"

# Update October 2022 (row 10) value
$ws.Range("I10").Value = 120

# Update Total 2022 (row 14) value
$ws.Range("I14").Value = 1255
